# "Generate Report for Handback" - update localization status report
# after a handback event for de-de (and a retroactive date fix for zh-cn).

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11500be126b55b93ea434b21616addfb63bc4221/e2e/"
$hyperlinkColor = 15570276   # OLE BGR value for RGB FF6495ED (the workbook's custom HyperLink font colour)

# ---------------------------------------------------------------------------
# 1) Overview sheet: status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both zh-cn and de-de columns,
#    and columns E/F grow wider to fit the new text.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"
$overview.Columns("E").ColumnWidth = 29.15
$overview.Columns("F").ColumnWidth = 29.15

# ---------------------------------------------------------------------------
# Helper: populate the "Latest Target File" (I) and "Latest Handback File"
# (J) columns for a locale sheet's data row, including the hyperlink and
# visual styling used by the other file-name hyperlinks in the workbook.
# ---------------------------------------------------------------------------
function Set-TargetFileLink {
    param($ws, [string]$cellRef, [string]$fileName)

    $ws.Hyperlinks.Add($ws.Range($cellRef), ($baseUrl + $fileName), [System.Type]::Missing, [System.Type]::Missing, $fileName) | Out-Null
    $ws.Range($cellRef).Font.Underline = 2
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: fill in Latest Target File / Latest Handback File, and fix
#    the Latest Handback DateTime (was the zero date, now the real one).
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

Set-TargetFileLink $zhcn "I2" "bc701561-146f-45ae-abb9-379e7256c3bc.md"
$zhcn.Range("J2").Value = "bc701561-146f-45ae-abb9-379e7256c3bc.6a9e2792bf47ca84268628c442d131b17944b134.zh-cn.xlf"

Set-TargetFileLink $zhcn "I3" "ccc4e32e-1cba-4d6e-8217-67b15d70c698.md"
$zhcn.Range("J3").Value = "ccc4e32e-1cba-4d6e-8217-67b15d70c698.cd01d35297dfb079af1e15ef009ef5bda9a33829.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-08-31 06:06:56"
$zhcn.Range("K3").Value = "2016-08-31 06:06:56"

$zhcn.Columns("C").ColumnWidth = 29.15
$zhcn.Columns("I").ColumnWidth = 39.15
$zhcn.Columns("J").ColumnWidth = 39.15

# ---------------------------------------------------------------------------
# 3) de-de sheet: same Latest Target File / Latest Handback File population,
#    plus the handback just completed so the Handback DateTime is newly set.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

Set-TargetFileLink $dede "I2" "bc701561-146f-45ae-abb9-379e7256c3bc.md"
$dede.Range("J2").Value = "bc701561-146f-45ae-abb9-379e7256c3bc.6a9e2792bf47ca84268628c442d131b17944b134.de-de.xlf"
$dede.Range("K2").Value = "2016-08-31 06:07:18"

Set-TargetFileLink $dede "I3" "ccc4e32e-1cba-4d6e-8217-67b15d70c698.md"
$dede.Range("J3").Value = "ccc4e32e-1cba-4d6e-8217-67b15d70c698.cd01d35297dfb079af1e15ef009ef5bda9a33829.de-de.xlf"
$dede.Range("K3").Value = "2016-08-31 06:07:18"

$dede.Columns("C").ColumnWidth = 29.15
$dede.Columns("I").ColumnWidth = 39.15
$dede.Columns("J").ColumnWidth = 39.15
